$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6.928617000579834
$ws.Range("B1").Value = 5.383649826049805
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 4.56832218170166
$ws.Range("E1").Value = 2.198020696640015
